# daily auto push: 2026-02-15 22:42 UTC
#
# The daily data-collection run appended one more hourly sample
# ("2026/02/16 月 4:00, rank 201") to the tail of the "2026/02/16" block
# that already starts at row 827. That pushes every later row down by
# one: the old row 828 ("2026/12/29 ...") becomes row 829, and so on,
# all the way to the old last row 869, which becomes row 870. The sheet
# dimension grows from A1:D869 to A1:D870 accordingly.
#
# We reproduce that with a genuine row insert (which shifts rows
# 828:869 down to 829:870 for us) and then populate the freshly
# inserted row 828 with the new sample.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 828 downward by inserting a new blank row at 828.
$ws.Rows(828).Insert()

# A828/B828 ("2026/02/16" / "月") are identical text to A827/B827, so
# copy them instead of re-typing the literal strings — a plain
# Range.Value assignment of a date-shaped string like "2026/02/16"
# gets auto-converted to a date serial by the COM layer, while copying
# an existing text cell preserves its literal text type untouched.
$ws.Range("A827:B827").Copy($ws.Range("A828"))

# C828/D828 are plain numbers, so a direct Value assignment is fine.
$ws.Range("C828").Value = 4
$ws.Range("D828").Value = 201
